$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with the latest scraped values.
# Numeric-looking price strings are forced to text (matching the original
# inline-string cell type) using a leading apostrophe, then the cell style is
# reset to "Normal" so no extra text-number-format style gets attached.

$ws.Range("D2").Value = '34.454.00'
$ws.Range("E2").Value = '  +0.95%  '
$ws.Range("D3").Value = '1.798.09'
$ws.Range("E3").Value = '  +0.55%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'227.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("E6").Value = '  +1.57%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").Value = "'32.46"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.78%  '
$ws.Range("D9").Value = "'0.297"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.51%  '
$ws.Range("D10").Value = "'0.0695"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.98%  '
$ws.Range("D11").Value = "'0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.60%  '
$ws.Range("D12").Value = '2.055.15'
$ws.Range("E12").Value = '  +0.41%  '
$ws.Range("D13").Value = '1.813.00'
$ws.Range("E13").Value = '  +1.37%  '
$ws.Range("D14").Value = "'11.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.89%  '
$ws.Range("D15").Value = "'0.636"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.56%  '
$ws.Range("D16").Value = '34.432.91'
$ws.Range("E16").Value = '  +1.10%  '
$ws.Range("D17").Value = "'4.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.77%  '
$ws.Range("D18").Value = "'68.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.71%  '
$ws.Range("D19").Value = "'247.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.63%  '
$ws.Range("D20").Value = '0.0₃0802'
$ws.Range("E20").Value = '  +3.22%  '
$ws.Range("D21").Value = "'11.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.66%  '
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").Value = "'4.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.58%  '
$ws.Range("D24").Value = "'2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").Value = "'163.18"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.97%  '
$ws.Range("E26").Value = '  +1.59%  '
$ws.Range("D27").Value = "'16.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.55%  '
$ws.Range("E28").Value = '  +2.53%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  +0.27%  '
$ws.Range("D31").Value = "'0.0522"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.86%  '
$ws.Range("E32").Value = '  +8.21%  '
$ws.Range("D33").Value = "'3.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.27%  '
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = '1.444.53'
$ws.Range("E35").Value = '  -1.20%  '
$ws.Range("D36").Value = "'2.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +6.94%  '
$ws.Range("E37").Value = '  +3.72%  '
$ws.Range("E38").Value = '  -0.26%  '
$ws.Range("D39").Value = "'1.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.80%  '
$ws.Range("D40").Value = "'84.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +4.79%  '
$ws.Range("E41").Value = '  +1.30%  '
$ws.Range("D42").Value = "'0.934"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.55%  '
$ws.Range("E43").Value = '  +2.60%  '
$ws.Range("D44").Value = "'13.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.16%  '
$ws.Range("E45").Value = '  +3.23%  '
$ws.Range("E46").Value = '  +0.62%  '
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").Value = '1.950.98'
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("D49").Value = "'105.77"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.41%  '
$ws.Range("E50").Value = '  -0.10%  '
$ws.Range("E51").Value = '  -4.46%  '
